# "Generate model answer with AI"
#
# The extracted student-answer image references were stored with a leading
# "/content" prefix (an artifact of the Colab/AI environment used to produce
# them). Normalize them to be relative paths so the workbook resolves the
# images correctly regardless of the runtime's working directory.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StudentAnswerExtracted")

$ws.Range("A2").Value = "images/WhatsApp Image 2024-04-28 at 4.27.06 PM.jpeg"
$ws.Range("A3").Value = "images/WhatsApp Image 2024-04-28 at 4.21.33 PM.jpeg"
$ws.Range("A4").Value = "images/WhatsApp Image 2024-04-28 at 4.20.17 PM.jpeg"
$ws.Range("A5").Value = "images/WhatsApp Image 2024-04-28 at 4.18.53 PM.jpeg"
$ws.Range("A6").Value = "images/WhatsApp Image 2024-04-28 at 4.17.26 PM.jpeg"

# Leave the cursor where the author last looked while reviewing the sheet.
$ws.Activate()
$ws.Range("E18").Select()
